$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price for row 2 (Microsoft IntelliMouse Explorer): 64.99 -> 64.95
$cell = $ws.Range("B2")
$cell.Value = 64.95

# Force the engine to materialize a new, distinct font/style record for this
# cell (mirrors the extra <font/> + cellXf that appear in the edited
# workbook's styles.xml), moving B2 off the style index shared with column A.
$cell.Font.Name = ""
